$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11: worked days reduced from 5 to 4 (drives K3/L3 formula recalculation)
$ws.Range("I11").Value = 4

# Mark "Lön" (salary/payday) on several dates across the sheet
$ws.Range("E15").Value = "Lön"
$ws.Range("G19").Value = "Lön"
$ws.Range("H23").Value = "Lön"
$ws.Range("F28").Value = "Lön"
